# Scheduled runner update: refresh Universalis market-price snapshots
# across the Leve profit sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 482.33334
$ws.Range("I28").Value = 482.33334
$ws.Range("K28").Value = 482.33334
$ws.Range("M28").Value = 2.666659999999979

$ws.Range("H98").Value = 20601.838
$ws.Range("J98").Value = 6508.6665
$ws.Range("L98").Value = 6508.6665
$ws.Range("N98").Value = -9504.666499999999

$ws.Range("H101").Value = 11913255
$ws.Range("I101").Value = 20415870
$ws.Range("J101").Value = 9595.6
$ws.Range("K101").Value = 61247610
$ws.Range("L101").Value = 28786.8
$ws.Range("M101").Value = -61245988
$ws.Range("N101").Value = -32030.8

$ws.Range("H107").Value = 5985.6313
$ws.Range("I107").Value = 7112.75
$ws.Range("K107").Value = 7112.75
$ws.Range("M107").Value = -5192.75

$ws.Range("H122").Value = 20601.838
$ws.Range("J122").Value = 6508.6665
$ws.Range("L122").Value = 19525.9995
$ws.Range("N122").Value = -24425.9995

$ws.Range("H137").Value = 6336.1167
$ws.Range("I137").Value = 8281.333000000001
$ws.Range("J137").Value = 2723.5715
$ws.Range("K137").Value = 24843.999
$ws.Range("L137").Value = 8170.7145
$ws.Range("M137").Value = -22293.999
$ws.Range("N137").Value = -13270.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2321.926
$ws.Range("I2").Value = 2293.68
$ws.Range("K2").Value = 2293.68
$ws.Range("M2").Value = -2180.68

$ws.Range("H32").Value = 5098.7
$ws.Range("I32").Value = 4967.0527
$ws.Range("K32").Value = 4967.0527
$ws.Range("M32").Value = -4680.0527

$ws.Range("H43").Value = 5342
$ws.Range("I43").Value = 5342
$ws.Range("K43").Value = 5342
$ws.Range("M43").Value = -5029

$ws.Range("H45").Value = 75849.64
$ws.Range("I45").Value = 94462.55
$ws.Range("K45").Value = 94462.55
$ws.Range("M45").Value = -94085.55

$ws.Range("H74").Value = 1518.093
$ws.Range("I74").Value = 828.8148
$ws.Range("K74").Value = 828.8148
$ws.Range("M74").Value = 45.18520000000001

$ws.Range("H77").Value = 1518.093
$ws.Range("I77").Value = 828.8148
$ws.Range("K77").Value = 4144.074
$ws.Range("M77").Value = 223.9260000000004

$ws.Range("H97").Value = 13341824
$ws.Range("I97").Value = 15608.286
$ws.Range("K97").Value = 15608.286
$ws.Range("M97").Value = -15112.286

$ws.Range("H102").Value = 13160.814
$ws.Range("I102").Value = 14584.435
$ws.Range("J102").Value = 4975
$ws.Range("K102").Value = 14584.435
$ws.Range("L102").Value = 4975
$ws.Range("M102").Value = -12962.435
$ws.Range("N102").Value = -8219

$ws.Range("H116").Value = 2321.926
$ws.Range("I116").Value = 2293.68
$ws.Range("K116").Value = 2293.68
$ws.Range("M116").Value = 0.3200000000001637

$ws.Range("H122").Value = 861610.7
$ws.Range("I122").Value = 3870.6072
$ws.Range("J122").Value = 4292571
$ws.Range("K122").Value = 11611.8216
$ws.Range("L122").Value = 12877713
$ws.Range("M122").Value = -9161.821599999999
$ws.Range("N122").Value = -12882613

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2321.926
$ws.Range("I3").Value = 2293.68
$ws.Range("K3").Value = 2293.68
$ws.Range("M3").Value = -2179.68

$ws.Range("H107").Value = 832.4091
$ws.Range("I107").Value = 711.35297
$ws.Range("K107").Value = 711.35297
$ws.Range("M107").Value = 1208.64703

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6173.4346
$ws.Range("I86").Value = 5128.294
$ws.Range("J86").Value = 9134.666999999999
$ws.Range("K86").Value = 5128.294
$ws.Range("L86").Value = 9134.666999999999
$ws.Range("M86").Value = -4005.294
$ws.Range("N86").Value = -11380.667

$ws.Range("H89").Value = 6173.4346
$ws.Range("I89").Value = 5128.294
$ws.Range("J89").Value = 9134.666999999999
$ws.Range("K89").Value = 25641.47
$ws.Range("L89").Value = 45673.335
$ws.Range("M89").Value = -20025.47
$ws.Range("N89").Value = -56905.335

$ws.Range("H99").Value = 4152819
$ws.Range("I99").Value = 7747416
$ws.Range("K99").Value = 7747416
$ws.Range("M99").Value = -7745918

$ws.Range("H105").Value = 102526.38
$ws.Range("I105").Value = 142170.94
$ws.Range("J105").Value = 3415
$ws.Range("K105").Value = 142170.94
$ws.Range("L105").Value = 3415
$ws.Range("M105").Value = -140423.94
$ws.Range("N105").Value = -6909

$ws.Range("H107").Value = 52638290
$ws.Range("I107").Value = 83343064
$ws.Range("J107").Value = 1528.1428
$ws.Range("K107").Value = 83343064
$ws.Range("L107").Value = 1528.1428
$ws.Range("M107").Value = -83341144
$ws.Range("N107").Value = -5368.1428

$ws.Range("H122").Value = 13751.4
$ws.Range("I122").Value = 16539.375
$ws.Range("K122").Value = 49618.125
$ws.Range("M122").Value = -47168.125

$ws.Range("H126").Value = 4152819
$ws.Range("I126").Value = 7747416
$ws.Range("K126").Value = 23242248
$ws.Range("M126").Value = -23239778

$ws.Range("H132").Value = 2197
$ws.Range("I132").Value = 1801.3334
$ws.Range("J132").Value = 3087.25
$ws.Range("K132").Value = 5404.0002
$ws.Range("L132").Value = 9261.75
$ws.Range("M132").Value = -2874.0002
$ws.Range("N132").Value = -14321.75

$ws.Range("H134").Value = 3130.2856
$ws.Range("I134").Value = 1484.1666
$ws.Range("K134").Value = 4452.4998
$ws.Range("M134").Value = -1917.4998

$ws.Range("H141").Value = 123287.95
$ws.Range("J141").Value = 127129.9
$ws.Range("L141").Value = 127129.9
$ws.Range("N141").Value = -137489.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 36210.668
$ws.Range("J136").Value = 36210.668
$ws.Range("L136").Value = 108632.004
$ws.Range("N136").Value = -113732.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5260.5
$ws.Range("I100").Value = 1770.5
$ws.Range("J100").Value = 8750.5
$ws.Range("K100").Value = 1770.5
$ws.Range("L100").Value = 8750.5
$ws.Range("M100").Value = -1229.5

$ws.Range("H136").Value = 4624.256
$ws.Range("I136").Value = 3628.6191
$ws.Range("J136").Value = 5574.636
$ws.Range("K136").Value = 10885.8573
$ws.Range("L136").Value = 16723.908
$ws.Range("M136").Value = -8335.8573
$ws.Range("N136").Value = -21823.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 250014990
$ws.Range("I2").Value = 250014990
$ws.Range("K2").Value = 250014990
$ws.Range("M2").Value = -250014878

$ws.Range("H107").Value = 17921.555
$ws.Range("I107").Value = 3577.8
$ws.Range("J107").Value = 23438.385
$ws.Range("K107").Value = 10733.4
$ws.Range("L107").Value = 70315.155
$ws.Range("M107").Value = -8813.400000000001
$ws.Range("N107").Value = -74155.155

$ws.Range("H126").Value = 29418.723
$ws.Range("I126").Value = 44472.637
$ws.Range("J126").Value = 5762.5713
$ws.Range("K126").Value = 133417.911
$ws.Range("L126").Value = 17287.7139
$ws.Range("M126").Value = -130947.911
$ws.Range("N126").Value = -22227.7139

